$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hour hand coordinates (columns B,C driven by column A): radius changed 33 -> 16,
# and the C-column (y) formula sign flipped from + to - in front of the SIN term.
$ws.Range("B2:B49").Formula = "=ROUND(64 + 16 * COS((PI() / 2) - A2 * (PI() /24)), 0)"
$ws.Range("C2:C49").Formula = "=ROUND(105 - 16 * SIN((PI() / 2) - A2 * (PI() / 24)),0)"

# Minute hand coordinates (columns F,G driven by column E): F (x) stays the same,
# G (y) formula sign flipped from + to - in front of the SIN term.
$ws.Range("G2:G61").Formula = "=ROUND(105 - 33 * SIN((PI() / 2) - E2 * (PI() / 30)),0)"

# Update the view: select H59 (this also clears the old scrolled topLeftCell state)
$ws.Range("H59").Select()
